# "workinng on template sheet" - update client/period, re-style the
# template header fields (bold labels + underline-style value cells),
# rework the data rows, and drop the last data row / shift the total up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header / template fields
# ---------------------------------------------------------------------
$ws.Range("D3").Value = "MC  Donalds"
$ws.Range("D4").Value = "Jun 30 2022"

# Bold-black labels (no fill) for the left-hand template captions.
$labelCells = @("C3", "C4", "C5")
foreach ($ref in $labelCells) {
    $r = $ws.Range($ref)
    $r.Font.Bold = $true
    $r.Font.Color = 0
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
    $r.WrapText = $true
}

# Value / fill-in cells get a thin bottom border ("underline") to mark
# where the answer goes - matches D3/D4/D5/L3/L4/O3/O4/E5.
$underlineCells = @("D3", "D4", "D5", "L3", "L4", "O3", "O4", "E5")
foreach ($ref in $underlineCells) {
    $r = $ws.Range($ref)
    $b = $r.Borders.Item(9)
    $b.Color = 0
    $b.Weight = 2
    $b.LineStyle = 1
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
    $r.WrapText = $true
}

# D5:E5 becomes a single merged "Bank Confirmation Control Sheet" cell.
$ws.Range("D5:E5").Merge()

# ---------------------------------------------------------------------
# 2) Data rows - values rewritten (rows 9-11), row 12 removed entirely
#    (old SUMMIT BANK row), total row shifts from 13 -> 12.
# ---------------------------------------------------------------------

# Row 9: MCB BANK / SAVING / PKR / Clifton
$ws.Range("B9").Value = 6
$ws.Range("D9").Value = 10201
$ws.Range("F9").Value = "PKR"
$ws.Range("H9").Value = 5000
$ws.Range("I9").Value = 2000
$ws.Range("J9").Value = 3000
$ws.Range("K9").Value = -1000
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = ""
$ws.Range("O9").Value = ""

# Row 10: MEEZAN BANK / SAVING / $ / Saddar
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = "MEEZAN BANK"
$ws.Range("D10").Value = 10205
$ws.Range("E10").Value = "SAVING"
$ws.Range("F10").Value = "$"
$ws.Range("G10").Value = "Saddar"
$ws.Range("H10").Value = -600
$ws.Range("I10").Value = 600
$ws.Range("J10").Value = -6000
$ws.Range("K10").Value = 6600

# Row 11: ALBARAKA / CURRENT / USD / Kemari
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = "ALBARAKA"
$ws.Range("D11").Value = 10206
$ws.Range("E11").Value = "CURRENT"
$ws.Range("F11").Value = "USD"
$ws.Range("G11").Value = "Kemari"
$ws.Range("H11").Value = 9000
$ws.Range("I11").Value = 5000
$ws.Range("J11").Value = 4000
$ws.Range("K11").Value = 1000

# Row 12 (old SUMMIT BANK entry) is removed entirely; the grand-total
# row (old row 13) shifts up into row 12.
$ws.Rows.Item(12).Delete()

$ws.Range("H12").Value = 13400

# ---------------------------------------------------------------------
# 3) Number format: switch the accounting format to the explicit
#    "_(* #,##0.00_)..." pattern used across the template.
# ---------------------------------------------------------------------
$acctFormat = '_(* #,##0.00_);_(* \(#,##0.00\);_(* "-"??_);_(@_)'
$ws.Range("H9:K12").NumberFormat = $acctFormat
$ws.Range("K3").NumberFormat = $acctFormat
$ws.Range("K4").NumberFormat = $acctFormat
$ws.Range("H7:K7").NumberFormat = $acctFormat

# ---------------------------------------------------------------------
# 4) Keep the used-range anchored at A1 (matches the template's stored
#    dimension) and leave the active selection on the new total cell.
# ---------------------------------------------------------------------
$ws.Range("A1").NumberFormat = "General"
$ws.Range("H12").Select()
